$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-5, columns E..T (Adm2-Calcrl NATMI TPM update)
$data = @{
    2 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.28161
        H = 0.84483
        M = 91.60947133333333
        N = 274.828414
        O = 0.7121576185577153
        P = 0.7121576185577152
        Q = 25.79814322218
        R = 232.18328899962
        S = 0.7121576185577153
        T = 0.7121576185577152
    }
    3 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.28161
        H = 0.84483
        M = 22.83185066666667
        N = 68.495552
        O = 0.1774912152792038
        P = 0.1774912152792038
        Q = 6.429677466239999
        R = 57.86709719616
        S = 0.1774912152792038
        T = 0.1774912152792038
    }
    4 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.28161
        H = 0.84483
        M = 7.077809999999999
        N = 21.23343
        O = 0.05502178149094856
        P = 0.05502178149094855
        Q = 1.9931820741
        R = 17.9386386669
        S = 0.05502178149094856
        T = 0.05502178149094855
    }
    5 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.28161
        H = 0.84483
        M = 7.117379
        N = 21.352137
        O = 0.05532938467213248
        P = 0.05532938467213247
        Q = 2.00432510019
        R = 18.03892590171
        S = 0.05532938467213248
        T = 0.05532938467213247
    }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
